$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 383, pushing the existing rows 383:436 down to 385:438.
$ws.Range("A383:A384").EntireRow.Insert()

# Populate new row 383 with its data.
$ws.Cells.Item(383, 1).Value2  = 6
$ws.Cells.Item(383, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(383, 3).Value2  = "Metropolitana"
$ws.Cells.Item(383, 4).Value2  = 45154
$ws.Cells.Item(383, 5).Value2  = 13
$ws.Cells.Item(383, 6).Value2  = 100112026
$ws.Cells.Item(383, 7).Value2  = "Haba"
$ws.Cells.Item(383, 8).Value2  = "Sin especificar"
$ws.Cells.Item(383, 9).Value2  = "Primera"
$ws.Cells.Item(383, 10).Value2 = 360
$ws.Cells.Item(383, 11).Value2 = 10000
$ws.Cells.Item(383, 12).Value2 = 11000
$ws.Cells.Item(383, 13).Value2 = 10361
$ws.Cells.Item(383, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(383, 15).Value2 = "Provincia de Copiapó"
$ws.Cells.Item(383, 16).Value2 = 414
$ws.Cells.Item(383, 17).Value2 = 25
$ws.Cells.Item(383, 18).Value2 = "Hortaliza"

# Populate new row 384 with its data.
$ws.Cells.Item(384, 1).Value2  = 6
$ws.Cells.Item(384, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(384, 3).Value2  = "Metropolitana"
$ws.Cells.Item(384, 4).Value2  = 45154
$ws.Cells.Item(384, 5).Value2  = 13
$ws.Cells.Item(384, 6).Value2  = 100112026
$ws.Cells.Item(384, 7).Value2  = "Haba"
$ws.Cells.Item(384, 8).Value2  = "Sin especificar"
$ws.Cells.Item(384, 9).Value2  = "Primera"
$ws.Cells.Item(384, 10).Value2 = 700
$ws.Cells.Item(384, 11).Value2 = 10000
$ws.Cells.Item(384, 12).Value2 = 12000
$ws.Cells.Item(384, 13).Value2 = 10714
$ws.Cells.Item(384, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(384, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(384, 16).Value2 = 429
$ws.Cells.Item(384, 17).Value2 = 25
$ws.Cells.Item(384, 18).Value2 = "Hortaliza"
